# "add area to Q files stn5"
# Adds an Area (cross-sectional area per segment) calculation alongside the
# existing Q (discharge) calculation on Sheet1, plus a small two-column
# summary (Atotal / Qtotal) next to it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New headers --------------------------------------------------------
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# --- Column G: per-segment area -----------------------------------------
# Row 2 is the "x=0" baseline segment, so it measures from 0 instead of the
# row above.
$ws.Range("G2").Formula = "=(D2-0)*B2/100"

# Row 3 stands alone (its neighbour above, D2, is not part of the D3:D9
# shared-formula fill below).
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# Rows 4-15 share one relative formula (fills all the way past the last
# data row down to row 15, matching the original sheet's pattern of
# pre-extending helper columns below the data).
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# --- Column H: running total of the new Area column ----------------------
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# --- Small summary block: Atotal / Qtotal side references ---------------
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- Re-fill D3:D9 as one contiguous relative formula --------------------
# (Originally each of D3..D9 carried its own independently-typed copy of
# this formula; re-entering it as a single fill turns it into one shared
# formula group, same as the rest of the sheet's helper columns.)
$ws.Range("D3:D9").Formula = "=(A3/100+(A4/100-A3/100)/2)"

# --- Match the author's final selection ----------------------------------
$ws.Range("J2:K2").Select()
